# Update "Elapsed Duration(Hrs)" values (column G) on several sheets to
# reflect the later PCM snapshot timestamps referenced in the commit.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "R1"; Cell = "G2"; Value = "3951:43:13" },
    @{ Sheet = "R1"; Cell = "G3"; Value = "91:15:51" },
    @{ Sheet = "R1"; Cell = "G4"; Value = "114:15:51" },
    @{ Sheet = "R2"; Cell = "G2"; Value = "12133:06:53" },
    @{ Sheet = "R2"; Cell = "G3"; Value = "3262:50:22" },
    @{ Sheet = "R2"; Cell = "G4"; Value = "501:01:56" },
    @{ Sheet = "R4"; Cell = "G2"; Value = "2978:56:42" },
    @{ Sheet = "R4"; Cell = "G3"; Value = "206:08:57" },
    @{ Sheet = "R4"; Cell = "G4"; Value = "94:21:22" },
    @{ Sheet = "R4"; Cell = "G5"; Value = "91:58:55" },
    @{ Sheet = "R5"; Cell = "G2"; Value = "452:55:41" },
    @{ Sheet = "R6"; Cell = "G2"; Value = "93:27:59" }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range($u.Cell).Value = $u.Value
}
